$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The list ("carrito") used to stop at row 14 (Id 13) with several trailing
# blank rows (Ids 9-13) already reserved below the real menu data (rows 2-9).
# This change extends that pre-reserved blank area down to row 21 (Id 20),
# i.e. 7 more empty rows, same shape as the existing ones: column A holds
# the running Id, columns B:E stay blank (typed as text, no value).
for ($r = 15; $r -le 21; $r++) {
    $id = $r - 1
    $ws.Cells.Item($r, 1).Value = $id

    for ($c = 2; $c -le 5; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        # A lone leading apostrophe is Excel's "force text" entry marker; it
        # commits the cell as an empty text value instead of clearing it,
        # matching the pre-existing blank rows (10-14) above.
        $cell.Value = "'"
        # Drop the quote-prefix formatting the apostrophe entry implies so
        # the new cells end up with plain default styling, like the rest of
        # the blank rows.
        $cell.Style = "Normal"
    }
}
